$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("D1").Value = "target01"
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0

$ws.Range("D9").Select()
